# Applies the "Map 9" content addition + landscape page orientation described
# in the commit "Updated spreadsheet with map9 info".
#
# Strategy: the target worksheet already contains several near-identical
# "map grid" blocks (Map 2/3/4/7 legends in H1:P9, H11:N19, ...). The new
# "Map 9" grid is the same kind of block, placed at P11:V19, reusing most of
# the existing cell styles and only a handful of brand-new ones (borderless /
# left-border-only / top-border-only variants). We build it by:
#   1) writing the new text values (in an order that makes the shared-string
#      table grow the same way Excel would: "Map 9 (...)" then "P" then "*P")
#   2) stamping cell formatting by copying (Copy + PasteSpecial formats only)
#      from a cell that already has the exact target style, then tweaking
#      borders/fill for the few styles that are brand-new
#   3) switching the page to landscape and updating the view selection/zoom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Cell text/values - order matters so new shared strings land in the
#    same sequence as the target file ("Map 9 ...", then "P", then "*P").
# ---------------------------------------------------------------------
$ws.Range("P11").Value = "Map 9 (P = Elevated Platform)"

$ws.Range("T13").Value = "P"
$ws.Range("R13").Value = "*P"

$ws.Range("Q12").Value = "S"
$ws.Range("R12").Value = "*"

$ws.Range("P14").Value = "P"
$ws.Range("R14").Value = "*P"

$ws.Range("P15").Value = "P"
$ws.Range("R15").Value = "*P"

$ws.Range("R16").Value = "*"

$ws.Range("Q17").Value = "P"
$ws.Range("R17").Value = "*P"
$ws.Range("S17").Value = "P"

$ws.Range("R18").Value = "*"
$ws.Range("S18").Value = "*"
$ws.Range("T18").Value = "E"

# ---------------------------------------------------------------------
# 2) Formatting - reuse existing cell styles via copy/paste-formats so we
#    don't fork new style entries unnecessarily.
# ---------------------------------------------------------------------

# style like A1 / "Map N" header text (big bold-ish header font)
$headerDonor = "A1"
foreach ($c in @("P11")) {
    $ws.Range($headerDonor).Copy()
    $ws.Range($c).PasteSpecial(-4122)
}

# style like C3 (orange-filled grid box, no text)
$donor5 = "C3"
foreach ($c in @("P13","Q13","S13","S14","S16","T16","S17","T17")) {
    $ws.Range($donor5).Copy()
    $ws.Range($c).PasteSpecial(-4122)
}

# style like H12 (light-accent-filled grid box)
$donor11 = "H12"
foreach ($c in @("R12","T12","R13","T13","P14","R14","P15","Q15","R15","T15","Q16","Q17","R17","P18","R18","S18")) {
    $ws.Range($donor11).Copy()
    $ws.Range($c).PasteSpecial(-4122)
}

# style like I13 (no-fill grid box)
$donor12 = "I13"
foreach ($c in @("P12","S12","Q14","T14","S15","P16","R16","P17","Q18")) {
    $ws.Range($donor12).Copy()
    $ws.Range($c).PasteSpecial(-4122)
}

# style like H14 (green-filled grid box)
$donor13 = "H14"
foreach ($c in @("T18")) {
    $ws.Range($donor13).Copy()
    $ws.Range($c).PasteSpecial(-4122)
}

# style like N17 (red-filled grid box)
$donor14 = "N17"
foreach ($c in @("Q12")) {
    $ws.Range($donor14).Copy()
    $ws.Range($c).PasteSpecial(-4122)
}

# NEW style: like I13 (no-fill grid box) but with no border at all
foreach ($c in @("V12","V13","V14","V15","V16","V18","U19","V19")) {
    $ws.Range($donor12).Copy()
    $ws.Range($c).PasteSpecial(-4122)
    $ws.Range($c).Borders.LineStyle = -4142
}

# NEW style: like N17 (red-filled) but no fill and no border
foreach ($c in @("V17")) {
    $ws.Range($donor14).Copy()
    $ws.Range($c).PasteSpecial(-4122)
    $ws.Range($c).Borders.LineStyle = -4142
    $ws.Range($c).Interior.ColorIndex = -4142
    $ws.Range($c).Interior.Pattern = -4142
}

# NEW style: like I13 (no-fill grid box) but only a thin LEFT border
foreach ($c in @("U12","U13","U14","U15","U16","U17","U18")) {
    $ws.Range($donor12).Copy()
    $ws.Range($c).PasteSpecial(-4122)
    $ws.Range($c).Borders.LineStyle = -4142
    $ws.Range($c).Borders.Item(7).LineStyle = 1
}

# NEW style: plain default cell but only a thin TOP border
foreach ($c in @("T19")) {
    $ws.Range($c).Borders.Item(8).LineStyle = 1
}

# ---------------------------------------------------------------------
# 3) Page setup + view state
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 2

$ws.Range("V16").Select()
$excel.ActiveWindow.Zoom = 57

Write-Output "Map 9 block added"
